$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns at D:E (existing D..K shift to F..M)
$ws.Range("D:E").EntireColumn.Insert()

# 2. Copy number formatting (date format) from F to the new D:E cells on the
#    three "Period Ending" header rows.
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("D80:E80").PasteSpecial(-4122)

# 3. Copy number formatting (#,##0) from F to the new D:E cells for every
#    other data row (contiguous blocks, skipping the blank separator rows).
$ws.Range("F8").Copy()
$ws.Range("D8:E35").PasteSpecial(-4122)
$ws.Range("D39:E77").PasteSpecial(-4122)
$ws.Range("D81:E102").PasteSpecial(-4122)

# 4. Populate the new column D (most recent quarter) and column E (prior
#    quarter) with the reported figures.
$ws.Range("D7").Value = 43464
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 128200
$ws.Range("E8").Value = 145300
$ws.Range("D9").Value = 107100
$ws.Range("E9").Value = 119600
$ws.Range("D10").Value = 21100
$ws.Range("E10").Value = 25700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 121400
$ws.Range("E17").Value = 134500
$ws.Range("D18").Value = 6800
$ws.Range("E18").Value = 10800
$ws.Range("D20").Value = -1200
$ws.Range("E20").Value = -700
$ws.Range("D21").Value = 11000
$ws.Range("E21").Value = 15600
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 5600
$ws.Range("E23").Value = 10100
$ws.Range("D24").Value = 1500
$ws.Range("E24").Value = 2700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 4100
$ws.Range("E26").Value = 7400
$ws.Range("D27").Value = 4100
$ws.Range("E27").Value = 7400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 1200
$ws.Range("E32").Value = 700
$ws.Range("D33").Value = 4100
$ws.Range("E33").Value = 7400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 4100
$ws.Range("E35").Value = 7400
$ws.Range("D38").Value = 43464
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 7300
$ws.Range("E41").Value = 4000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 61500
$ws.Range("E43").Value = 66400
$ws.Range("D44").Value = 67700
$ws.Range("E44").Value = 67800
$ws.Range("D45").Value = 5700
$ws.Range("E45").Value = 1700
$ws.Range("D46").Value = 142200
$ws.Range("E46").Value = 139900
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 117600
$ws.Range("E48").Value = 119400
$ws.Range("D49").Value = 125400
$ws.Range("E49").Value = 126800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 3200
$ws.Range("E52").Value = 3800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 388400
$ws.Range("E54").Value = 389800
$ws.Range("D57").Value = 29300
$ws.Range("E57").Value = 35500
$ws.Range("D58").Value = 9900
$ws.Range("E58").Value = 9900
$ws.Range("D59").Value = 13900
$ws.Range("E59").Value = 13500
$ws.Range("D60").Value = 53100
$ws.Range("E60").Value = 58800
$ws.Range("D61").Value = 80600
$ws.Range("E61").Value = 75800
$ws.Range("D62").Value = 37200
$ws.Range("E62").Value = 37400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 170900
$ws.Range("E66").Value = 172000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 163100
$ws.Range("E72").Value = 161400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 217500
$ws.Range("E76").Value = 217800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43464
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 4100
$ws.Range("E81").Value = 7400
$ws.Range("D83").Value = 5400
$ws.Range("E83").Value = 5500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 6200
$ws.Range("E89").Value = 18000
$ws.Range("D91").Value = -3000
$ws.Range("E91").Value = -1800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -2900
$ws.Range("E94").Value = -1800
$ws.Range("D96").Value = -2400
$ws.Range("E96").Value = -2400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = -17400
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 3300
$ws.Range("E102").Value = -1100

# 5. One genuine data correction that is not just a product of the column
#    shift: row 96 ("Dividends Paid") column J (old column H, "12/31/2016")
#    is restated from -4500 to 0.
$ws.Range("J96").Value = 0
